# "atualizei dados bibi e add" — update BIBI daily-revenue data:
#   1) June/2025 gets a new day-30 entry (24329.12) appended to the end of
#      its block, which currently runs rows 2-30 (days 1-29). Insert a row
#      right after the existing June rows (before May's block) and fill it.
#   2) The March/2025 block (31 rows, at the end of the original sheet) is
#      removed entirely, along with its now-unused "03/2025" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert the missing June 30th row right after row 30 (June 29th) ---
$ws.Rows.Item(31).Insert()
$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = 24329.12
$ws.Cells.Item(31, 3).Value = 6
$ws.Cells.Item(31, 4).Value = 2025
$ws.Cells.Item(31, 5).Value = "06/2025"

# --- 2) Delete the entire March/2025 block ---
# Original March block was rows 92:122; after inserting the row above,
# everything shifted down by one, so it now lives at rows 93:123.
$ws.Range("A93:E123").EntireRow.Delete()
